# Add a new worksheet "test14" after the last existing sheet ("test13"),
# mirroring the structure used by the other "testN" sheets (a 5-column by
# 8-row table of strings), but anchored at D68:H75 instead of A1:E8.
# Also updates the active tab so the new sheet becomes the selected one,
# which Excel does automatically by removing tabSelected from the
# previously active sheet (test13) and adding it to the new sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current last sheet (test13).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "test14"

# Copy the bold/filled header formatting used by the other "testN" sheets
# (row 1, columns A:E) onto the new header row (row 68, columns D:H) so
# the same cell style gets reused rather than creating a brand-new one.
$srcHeader = $wb.Worksheets.Item(2).Range("A1:E1")
$srcHeader.Copy()
$newSheet.Range("D68:H68").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Header row: a, b, c, d, e
$headers = @("a", "b", "c", "d", "e")
for ($c = 0; $c -lt 5; $c++) {
  $newSheet.Cells.Item(68, 4 + $c).Value = $headers[$c]
}

# Data rows mirror the other testN sheets (A1..E7 style values), except the
# "D" column of the 6th data row (row 74) carries the sheet's own unique
# marker string "blah14" instead of the generic "D6" value.
$data = @(
  @("A1", "B1", "C1", "D1", "E1"),
  @("A2", "B2", "C2", "D2", "E2"),
  @("A3", "B3", "C3", "D3", "E3"),
  @("A4", "B4", "C4", "D4", "E4"),
  @("A5", "B5", "C5", "D5", "E5"),
  @("A6", "B6", "C6", "blah14", "E6"),
  @("A7", "B7", "C7", "D7", "E7")
)

for ($r = 0; $r -lt 7; $r++) {
  for ($c = 0; $c -lt 5; $c++) {
    $newSheet.Cells.Item(69 + $r, 4 + $c).Value = $data[$r][$c]
  }
}

# Match the page setup used by the sibling sheets.
$pageSetup = $newSheet.PageSetup
$pageSetup.PaperSize = 9      # xlPaperA4
$pageSetup.Orientation = 1    # xlPortrait

# Leave the new sheet active with G74 selected, matching the saved view
# state captured in the workbook.
$newSheet.Activate()
$newSheet.Range("G74").Select()
